$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.935.24"
$ws.Range("E2").Value = "  -5.20%  "

$ws.Range("D3").Value = "3.605.41"
$ws.Range("E3").Value = "  -2.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.82%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "401.41"
$ws.Range("E5").Value = "  -5.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "124.47"
$ws.Range("E6").Value = "  -4.38%  "

$ws.Range("D7").Value = "3.606.88"
$ws.Range("E7").Value = "  -1.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  -7.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.998"
$ws.Range("E9").Value = "  -0.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.689"
$ws.Range("E10").Value = "  -9.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  -13.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000337"
$ws.Range("E12").Value = "  -5.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.30"
$ws.Range("E13").Value = "  -6.93%  "

$ws.Range("D14").Value = "4.144.13"
$ws.Range("E14").Value = "  -3.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.31"
$ws.Range("E15").Value = "  -6.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.39"
$ws.Range("E16").Value = "  +12.89%  "

$ws.Range("E17").Value = "  -2.56%  "

$ws.Range("D18").Value = "3.595.92"
$ws.Range("E18").Value = "  -2.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.87"
$ws.Range("E19").Value = "  -7.78%  "

$ws.Range("D20").Value = "63.979.54"
$ws.Range("E20").Value = "  -5.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.03"
$ws.Range("E21").Value = "  -8.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "398.20"
$ws.Range("E22").Value = "  -10.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.03"
$ws.Range("E23").Value = "  -8.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.57"
$ws.Range("E24").Value = "  -7.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.95"
$ws.Range("E25").Value = "  -4.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.43"
$ws.Range("E26").Value = "  +9.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "34.54"
$ws.Range("E27").Value = "  -7.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.04"
$ws.Range("E28").Value = "  -7.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.86"
$ws.Range("E29").Value = "  -14.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.05"
$ws.Range("E30").Value = "  -2.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("E31").Value = "  -3.24%  "

$ws.Range("E32").Value = "  -5.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.90"
$ws.Range("E33").Value = "  -3.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.150"
$ws.Range("E34").Value = "  -6.04%  "

$ws.Range("E35").Value = "  +0.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "37.05"
$ws.Range("E36").Value = "  -8.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.34"
$ws.Range("E37").Value = "  -3.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0440"
$ws.Range("E38").Value = "  -10.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.997"
$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0670"
$ws.Range("E40").Value = "  -6.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("E41").Value = "  -8.51%  "

$ws.Range("E42").Value = "  -9.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.08"
$ws.Range("E43").Value = "  +16.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "142.27"
$ws.Range("E44").Value = "  -2.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.12"
$ws.Range("E45").Value = "  -2.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.98"
$ws.Range("E46").Value = "  -4.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.11"
$ws.Range("E47").Value = "  -8.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.11"
$ws.Range("E48").Value = "  -4.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.49"
$ws.Range("E49").Value = "  -6.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.68"
$ws.Range("E50").Value = "  -8.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.281"
$ws.Range("E51").Value = "  -7.29%  "

